$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 12 ("Hello World!" demo slide)
# ---------------------------------------------------------------------------
$s12 = $p.Slides.Item(12)

# Title 1: "Hello World!\t" -> "When should I use it?" + "\t" (two runs)
$title = $s12.Shapes.Item(1).TextFrame.TextRange
$title.Text = "When should I use it?"
$title.InsertAfter("`t") | Out-Null

# Text Placeholder 2: replace single line with a 5-item bullet list, the
# 4th item containing a mid-sentence "etc" run.
$body = $s12.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Chat/Messaging Apps"
$body.InsertAfter("`rReal time Apps ( stocks / ticker tape)") | Out-Null
$body.InsertAfter("`rHighly Concurrent apps") | Out-Null
$body.InsertAfter("`rSingle page apps with lots of asynchronous calls (Gmail ") | Out-Null
$body.InsertAfter("etc") | Out-Null
$body.InsertAfter(")") | Out-Null
$body.InsertAfter("`rServing lots of dynamic content") | Out-Null

# ---------------------------------------------------------------------------
# Slide 4 (Timeline)
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)

# TextBox 18: merge the "EBay"/" releases API built on "/"node" runs down to
# "EBay" + " releases API built on node" (2 runs instead of 3).
$tb18 = $s4.Shapes.Item(7).TextFrame.TextRange
$ebayPara = $tb18.Paragraphs(3, 1)
$ebayRest = $ebayPara.Characters(5, $ebayPara.Length - 4)
$ebayRest.Text = " releases API built on node"

# TextBox 24: merge "Nov " + "2010  Cloud9IDE" into a single run.
$tb24 = $s4.Shapes.Item(19).TextFrame.TextRange
$novPara = $tb24.Paragraphs(1, 1)
$novAll = $novPara.Characters(1, $novPara.Length)
$novAll.Text = "Nov 2010  Cloud9IDE"
